$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.812.32"
$ws.Range("E2").Value = "  +5.30%  "

$ws.Range("D3").Value = "2.269.56"
$ws.Range("E3").Value = "  +3.63%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").Value = "301.75"
$ws.Range("E5").Value = "  +4.00%  "

$ws.Range("D6").Value = "92.59"
$ws.Range("E6").Value = "  +7.94%  "

$ws.Range("E7").Value = "  +3.66%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").Value = "0.487"
$ws.Range("E9").Value = "  +5.21%  "

$ws.Range("D10").Value = "54.66"
$ws.Range("E10").Value = "  +9.38%  "

$ws.Range("E11").Value = "  +8.49%  "

$ws.Range("D12").Value = "0.0799"
$ws.Range("E12").Value = "  +3.11%  "

$ws.Range("E13").Value = "  +3.08%  "

$ws.Range("D14").Value = "6.68"
$ws.Range("E14").Value = "  +4.49%  "

$ws.Range("D15").Value = "2.622.72"
$ws.Range("E15").Value = "  +3.85%  "

$ws.Range("D16").Value = "14.17"
$ws.Range("E16").Value = "  +3.91%  "

$ws.Range("D17").Value = "2.296.62"
$ws.Range("E17").Value = "  +3.42%  "

$ws.Range("D18").Value = "0.754"
$ws.Range("E18").Value = "  +4.27%  "

$ws.Range("D19").Value = "41.739.70"
$ws.Range("E19").Value = "  +5.42%  "

$ws.Range("D20").Value = "12.27"
$ws.Range("E20").Value = "  +10.52%  "

$ws.Range("D21").Value = "0.0₃0906"
$ws.Range("E21").Value = "  +3.13%  "

$ws.Range("E22").Value = "  +4.38%  "

$ws.Range("D23").Value = "67.22"
$ws.Range("E23").Value = "  +3.53%  "

$ws.Range("D24").Value = "'240.80"
$ws.Range("E24").Value = "  +2.06%  "

$ws.Range("D25").Value = "2.58"
$ws.Range("E25").Value = "  +5.89%  "

$ws.Range("E26").Value = "  -0.20%  "

$ws.Range("D27").Value = "1.89"
$ws.Range("E27").Value = "  +5.21%  "

$ws.Range("D28").Value = "23.79"
$ws.Range("E28").Value = "  +3.35%  "

$ws.Range("D29").Value = "'2.20"
$ws.Range("E29").Value = "  +2.70%  "

$ws.Range("D30").Value = "9.67"
$ws.Range("E30").Value = "  +5.87%  "

$ws.Range("D31").Value = "34.05"
$ws.Range("E31").Value = "  +9.53%  "

$ws.Range("D32").Value = "157.96"
$ws.Range("E32").Value = "  +1.55%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("E34").Value = "  +6.03%  "

$ws.Range("D35").Value = "0.0737"
$ws.Range("E35").Value = "  +4.92%  "

$ws.Range("D36").Value = "3.06"
$ws.Range("E36").Value = "  +8.53%  "

$ws.Range("E37").Value = "  +3.32%  "

$ws.Range("E38").Value = "  +7.49%  "

$ws.Range("D39").Value = "16.53"
$ws.Range("E39").Value = "  +9.86%  "

$ws.Range("E40").Value = "  +2.93%  "

$ws.Range("E41").Value = "  +7.16%  "

$ws.Range("D42").Value = "3.98"
$ws.Range("E42").Value = "  +7.74%  "

$ws.Range("D43").Value = "'20.20"
$ws.Range("E43").Value = "  +17.11%  "

$ws.Range("D44").Value = "2.054.41"
$ws.Range("E44").Value = "  -2.72%  "

$ws.Range("E45").Value = "  +4.55%  "

$ws.Range("D46").Value = "10.12"
$ws.Range("E46").Value = "  +4.67%  "

$ws.Range("E47").Value = "  +10.98%  "

$ws.Range("E48").Value = "  -3.98%  "

$ws.Range("D49").Value = "2.495.90"
$ws.Range("E49").Value = "  +4.28%  "

$ws.Range("E50").Value = "  +3.43%  "

$ws.Range("E51").Value = "  +5.24%  "
